# Update Google Drive "image" links in column C (rows 2-25) so that the
# uc?id= download links include the export=download parameter, i.e.
#   https://drive.google.com/uc?id=XXXX
# becomes
#   https://drive.google.com/uc?export=download&id=XXXX

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $val = $cell.Value2
    if ($val -ne $null -and $val -like "*drive.google.com/uc?id=*") {
        $newVal = $val -replace [regex]::Escape("uc?id="), "uc?export=download&id="
        $cell.Value2 = $newVal
    }
}

$wb.Save()
